$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated OSiRIX synth landmark measurements (column B values)
$ws.Range("B5").Value = 10278
$ws.Range("B6").Value = 8793
$ws.Range("B12").Value = 5335

# Match the active cell / selection left behind by the edit session
$ws.Range("O6").Select()
